# completed tank titrations 0330
# Append the new CRM accuracy reading taken on 2022-03-30 as row 74.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$ws.Range("A74").Value = 20220330
$ws.Range("B74").Value = 2227.7080000000001
$ws.Range("C74").Value = 2224.4699999999998
$ws.Range("D74").Formula = "=100*(B74-C74)/C74"
$ws.Range("E74").Value = 180
$ws.Range("F74").Value = "CRM OPENED 20220318"

# Match the author's final selection/viewport state in the saved workbook
$ws.Range("H72").Select()
